$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 10: Polygon -> OKB
Set-TextValue $ws.Range("B10") "OKB"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D10") "49.75"
Set-TextValue $ws.Range("E10") "  -2.23%  "

# Row 11: OKB -> Polygon
Set-TextValue $ws.Range("B11") "Polygon"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D11") "1.359"
Set-TextValue $ws.Range("E11") "  -3.05%  "

# Row 2
Set-TextValue $ws.Range("D2") "23.829.20"
Set-TextValue $ws.Range("E2") "  -3.39%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.616.79"
Set-TextValue $ws.Range("E3") "  -3.50%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "308.22"
Set-TextValue $ws.Range("E5") "  -1.86%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.04%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.3933"
Set-TextValue $ws.Range("E7") "  -0.43%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3835"
Set-TextValue $ws.Range("E8") "  -2.93%  "

# Row 9
Set-TextValue $ws.Range("D9") "1.001"
Set-TextValue $ws.Range("E9") "  -0.05%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.08440"
Set-TextValue $ws.Range("E12") "  -2.43%  "

# Row 13
Set-TextValue $ws.Range("D13") "23.79"
Set-TextValue $ws.Range("E13") "  -6.12%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.045"
Set-TextValue $ws.Range("E14") "  -3.80%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.585"
Set-TextValue $ws.Range("E15") "  -1.57%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.00001278"
Set-TextValue $ws.Range("E16") "  -3.22%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.614.02"
Set-TextValue $ws.Range("E17") "  -4.34%  "

# Row 18
Set-TextValue $ws.Range("D18") "93.75"
Set-TextValue $ws.Range("E18") "  -0.25%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06929"
Set-TextValue $ws.Range("E19") "  -1.09%  "

# Row 20
Set-TextValue $ws.Range("D20") "20.01"
Set-TextValue $ws.Range("E20") "  -5.25%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.808"
Set-TextValue $ws.Range("E21") "  -4.06%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.9995"
Set-TextValue $ws.Range("E22") "  -0.15%  "

# Row 23
Set-TextValue $ws.Range("D23") "13.40"
Set-TextValue $ws.Range("E23") "  -3.81%  "

# Row 24
Set-TextValue $ws.Range("D24") "23.835.94"
Set-TextValue $ws.Range("E24") "  -3.39%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.462"
Set-TextValue $ws.Range("E25") "  +5.08%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.848"
Set-TextValue $ws.Range("E26") "  +2.52%  "

# Row 27
Set-TextValue $ws.Range("D27") "22.21"
Set-TextValue $ws.Range("E27") "  -3.48%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -2.37%  "

# Row 29
Set-TextValue $ws.Range("D29") "139.83"
Set-TextValue $ws.Range("E29") "  -4.19%  "

# Row 30
Set-TextValue $ws.Range("D30") "5.278"
Set-TextValue $ws.Range("E30") "  -9.60%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.836"
Set-TextValue $ws.Range("E31") "  -5.89%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.502"
Set-TextValue $ws.Range("E32") "  +0.29%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.792.13"
Set-TextValue $ws.Range("E33") "  -3.78%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.08115"
Set-TextValue $ws.Range("E34") "  -1.84%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.9783"
Set-TextValue $ws.Range("E35") "  -1.58%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.02890"
Set-TextValue $ws.Range("E36") "  -6.26%  "

# Row 37
Set-TextValue $ws.Range("D37") "6.577"
Set-TextValue $ws.Range("E37") "  -5.60%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.2664"
Set-TextValue $ws.Range("E38") "  -5.35%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.09137"
Set-TextValue $ws.Range("E39") "  -5.02%  "

# Row 40
Set-TextValue $ws.Range("D40") "10.38"
Set-TextValue $ws.Range("E40") "  +0.68%  "

# Row 41
Set-TextValue $ws.Range("D41") "13.68"
Set-TextValue $ws.Range("E41") "  +0.95%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.427"
Set-TextValue $ws.Range("E42") "  -5.69%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.7499"
Set-TextValue $ws.Range("E43") "  -5.29%  "

# Row 44
Set-TextValue $ws.Range("D44") "16.04"
Set-TextValue $ws.Range("E44") "  -3.98%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.6896"
Set-TextValue $ws.Range("E45") "  -3.06%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.466"
Set-TextValue $ws.Range("E46") "  -3.99%  "

# Row 47
Set-TextValue $ws.Range("E47") "  -2.35%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.0000"
Set-TextValue $ws.Range("E48") "  -0.11%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.08221"
Set-TextValue $ws.Range("E49") "  -5.11%  "

# Row 50
Set-TextValue $ws.Range("D50") "134.16"
Set-TextValue $ws.Range("E50") "  -2.86%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.201"
Set-TextValue $ws.Range("E51") "  -9.60%  "
